$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string partial text edits (rich text runs)
$ws.Range("A8").Characters(21, 2).Text = "21"
$ws.Range("C9").Characters(27, 9).Text = "5/20/2024"
$ws.Range("C9").Characters(47, 9).Text = "5/26/2024"

# Numeric cell updates (crime statistics table)
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 66.666666666666
$ws.Range("F14").Value = 31
$ws.Range("G14").Value = 31
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 139
$ws.Range("J14").Value = 162
$ws.Range("K14").Value = -14.197530864197
$ws.Range("L14").Value = -21.468926553672
$ws.Range("M14").Value = -24.043715846994
$ws.Range("N14").Value = -81.317204301075
$ws.Range("C15").Value = 26
$ws.Range("D15").Value = 28
$ws.Range("E15").Value = -7.142857142857
$ws.Range("F15").Value = 120
$ws.Range("G15").Value = 105
$ws.Range("H15").Value = 14.285714285714
$ws.Range("I15").Value = 611
$ws.Range("J15").Value = 599
$ws.Range("K15").Value = 2.003338898163
$ws.Range("L15").Value = -4.828660436137
$ws.Range("M15").Value = 23.434343434343
$ws.Range("N15").Value = -51.73775671406
$ws.Range("C16").Value = 324
$ws.Range("D16").Value = 306
$ws.Range("E16").Value = 5.882352941176
$ws.Range("F16").Value = 1259
$ws.Range("G16").Value = 1223
$ws.Range("H16").Value = 2.943581357318
$ws.Range("I16").Value = 6489
$ws.Range("J16").Value = 6158
$ws.Range("K16").Value = 5.375121792789
$ws.Range("L16").Value = 2.463287541449
$ws.Range("M16").Value = -8.682803264846
$ws.Range("N16").Value = -80.361952607208
$ws.Range("C17").Value = 630
$ws.Range("D17").Value = 569
$ws.Range("E17").Value = 10.720562390158
$ws.Range("F17").Value = 2346
$ws.Range("G17").Value = 2176
$ws.Range("H17").Value = 7.8125
$ws.Range("I17").Value = 11014
$ws.Range("J17").Value = 10510
$ws.Range("K17").Value = 4.795432921027
$ws.Range("L17").Value = 13.745739956625
$ws.Range("M17").Value = 70.918684047175
$ws.Range("N17").Value = -28.121125106049
$ws.Range("C18").Value = 240
$ws.Range("D18").Value = 256
$ws.Range("E18").Value = -6.25
$ws.Range("F18").Value = 941
$ws.Range("G18").Value = 1038
$ws.Range("H18").Value = -9.344894026974
$ws.Range("I18").Value = 5073
$ws.Range("J18").Value = 5744
$ws.Range("K18").Value = -11.681754874651
$ws.Range("L18").Value = -16.781496062992
$ws.Range("M18").Value = -27.435273923616
$ws.Range("N18").Value = -87.045785347667
$ws.Range("C19").Value = 914
$ws.Range("D19").Value = 978
$ws.Range("E19").Value = -6.543967280163
$ws.Range("F19").Value = 3658
$ws.Range("G19").Value = 3853
$ws.Range("H19").Value = -5.060991435245
$ws.Range("I19").Value = 18872
$ws.Range("J19").Value = 19546
$ws.Range("K19").Value = -3.448275862068
$ws.Range("L19").Value = -3.36917562724
$ws.Range("M19").Value = 35.496840896036
$ws.Range("N19").Value = -41.025
$ws.Range("C20").Value = 285
$ws.Range("D20").Value = 304
$ws.Range("E20").Value = -6.25
$ws.Range("F20").Value = 1096
$ws.Range("G20").Value = 1232
$ws.Range("H20").Value = -11.038961038961
$ws.Range("I20").Value = 5309
$ws.Range("J20").Value = 5933
$ws.Range("K20").Value = -10.517444800269
$ws.Range("L20").Value = 4.323049715071
$ws.Range("M20").Value = 35.020345879959
$ws.Range("N20").Value = -88.214793109572
$ws.Range("C21").Value = 2429
$ws.Range("D21").Value = 2447
$ws.Range("E21").Value = -0.735594605639
$ws.Range("F21").Value = 9451
$ws.Range("G21").Value = 9658
$ws.Range("H21").Value = -2.143300890453
$ws.Range("I21").Value = 47507
$ws.Range("J21").Value = 48652
$ws.Range("K21").Value = -2.353448984625
$ws.Range("L21").Value = -0.090431125131
$ws.Range("M21").Value = 21.566570280713
$ws.Range("N21").Value = -71.481826094786
$ws.Range("C22").Value = 47
$ws.Range("D22").Value = 37
$ws.Range("E22").Value = 27.027027027027
$ws.Range("F22").Value = 166
$ws.Range("G22").Value = 174
$ws.Range("H22").Value = -4.597701149425
$ws.Range("I22").Value = 836
$ws.Range("J22").Value = 880
$ws.Range("K22").Value = -5
$ws.Range("L22").Value = -11.627906976744
$ws.Range("M22").Value = 1.456310679611
$ws.Range("C23").Value = 132
$ws.Range("D23").Value = 113
$ws.Range("E23").Value = 16.814159292035
$ws.Range("F23").Value = 467
$ws.Range("G23").Value = 467
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 2364
$ws.Range("J23").Value = 2440
$ws.Range("K23").Value = -3.11475409836
$ws.Range("L23").Value = 4.973357015985
$ws.Range("M23").Value = 56.660039761431
$ws.Range("C24").Value = 2165
$ws.Range("D24").Value = 2296
$ws.Range("E24").Value = -5.705574912891
$ws.Range("F24").Value = 8077
$ws.Range("G24").Value = 8496
$ws.Range("H24").Value = -4.931732580037
$ws.Range("I24").Value = 42542
$ws.Range("J24").Value = 42681
$ws.Range("K24").Value = -0.325671844614
$ws.Range("L24").Value = -1.589210946355
$ws.Range("M24").Value = 40.458267300581
$ws.Range("C25").Value = 1224
$ws.Range("D25").Value = 1282
$ws.Range("E25").Value = -4.524180967238
$ws.Range("F25").Value = 4487
$ws.Range("G25").Value = 4575
$ws.Range("H25").Value = -1.923497267759
$ws.Range("I25").Value = 23977
$ws.Range("J25").Value = 22972
$ws.Range("K25").Value = 4.374891171861
$ws.Range("L25").Value = -1.231669138243
$ws.Range("C26").Value = 1064
$ws.Range("D26").Value = 868
$ws.Range("E26").Value = 22.58064516129
$ws.Range("F26").Value = 3930
$ws.Range("G26").Value = 3621
$ws.Range("H26").Value = 8.533554266777
$ws.Range("I26").Value = 17987
$ws.Range("J26").Value = 16786
$ws.Range("K26").Value = 7.15477183367
$ws.Range("L26").Value = 12.997864053273
$ws.Range("M26").Value = 1.466689231116
$ws.Range("C27").Value = 39
$ws.Range("D27").Value = 58
$ws.Range("E27").Value = -32.758620689655
$ws.Range("F27").Value = 186
$ws.Range("G27").Value = 207
$ws.Range("H27").Value = -10.144927536231
$ws.Range("I27").Value = 972
$ws.Range("J27").Value = 997
$ws.Range("K27").Value = -2.507522567703
$ws.Range("L27").Value = -8.474576271186
$ws.Range("C28").Value = 112
$ws.Range("D28").Value = 113
$ws.Range("E28").Value = -0.884955752212
$ws.Range("F28").Value = 464
$ws.Range("G28").Value = 459
$ws.Range("H28").Value = 1.089324618736
$ws.Range("I28").Value = 2064
$ws.Range("J28").Value = 2047
$ws.Range("K28").Value = 0.830483634587
$ws.Range("L28").Value = 5.629477993858
$ws.Range("C29").Value = 28
$ws.Range("D29").Value = 15
$ws.Range("E29").Value = 86.666666666666
$ws.Range("F29").Value = 100
$ws.Range("G29").Value = 83
$ws.Range("H29").Value = 20.481927710843
$ws.Range("I29").Value = 375
$ws.Range("J29").Value = 428
$ws.Range("K29").Value = -12.383177570093
$ws.Range("L29").Value = -34.325744308231
$ws.Range("M29").Value = -36.868686868686
$ws.Range("N29").Value = -83.092876465284
$ws.Range("C30").Value = 24
$ws.Range("D30").Value = 14
$ws.Range("E30").Value = 71.428571428571
$ws.Range("F30").Value = 83
$ws.Range("G30").Value = 73
$ws.Range("H30").Value = 13.698630136986
$ws.Range("I30").Value = 317
$ws.Range("J30").Value = 363
$ws.Range("K30").Value = -12.672176308539
$ws.Range("L30").Value = -34.77366255144
$ws.Range("M30").Value = -34.907597535934
$ws.Range("N30").Value = -84.19740777667
$ws.Range("C31").Value = 9
$ws.Range("D31").Value = 19
$ws.Range("E31").Value = -52.631578947368
$ws.Range("F31").Value = 61
$ws.Range("G31").Value = 53
$ws.Range("H31").Value = 15.094339622641
$ws.Range("I31").Value = 254
$ws.Range("J31").Value = 228
$ws.Range("K31").Value = 11.403508771929
$ws.Range("L31").Value = -8.960573476702
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 3
$ws.Range("E33").Value = 0
$ws.Range("I33").Value = 93
$ws.Range("J33").Value = 91
$ws.Range("K33").Value = 2.197802197802
$ws.Range("L33").Value = 6.896551724137
